# feat: update preferences order evalution
#
# Updates the Mean (column H) and Std (column I) values for the
# PreferenceOrder / PartialOrder evaluation rows on the "Overall" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Overall")

$updates = @(
    @{Row=26; H=0.40971; I=0.01701},
    @{Row=27; H=0.01586; I=0.00996},
    @{Row=28; H=0.38551; I=0.0187},
    @{Row=29; H=0.00072; I=0.00244},
    @{Row=30; H=0.41069; I=0.01804},
    @{Row=31; H=0.01586; I=0.01089},
    @{Row=32; H=0.38124; I=0.01977},
    @{Row=33; H=0.00108; I=0.00293},
    @{Row=34; H=0.48324; I=0.01144},
    @{Row=35; H=0.00072; I=0.00244},
    @{Row=36; H=0.48404; I=0.01228},
    @{Row=37; H=0.00072; I=0.00244},
    @{Row=38; H=0.48341; I=0.01132},
    @{Row=39; H=0.00072; I=0.00244},
    @{Row=40; H=0.4836; I=0.01222},
    @{Row=41; H=0.00072; I=0.00244},
    @{Row=66; H=0.38145; I=0.01908},
    @{Row=67; H=0.01153; I=0.00865},
    @{Row=68; H=0.36819; I=0.01595},
    @{Row=69; H=0.00036; I=0.00177},
    @{Row=70; H=0.38441; I=0.0183},
    @{Row=71; H=0.01117; I=0.00856},
    @{Row=72; H=0.36718; I=0.01696},
    @{Row=73; H=0.00108; I=0.00293},
    @{Row=74; H=0.46369; I=0.01088},
    @{Row=75; H=0.00072; I=0.00244},
    @{Row=76; H=0.46825; I=0.01062},
    @{Row=77; H=0.00144; I=0.0033},
    @{Row=78; H=0.46395; I=0.01055},
    @{Row=79; H=0.00072; I=0.00244},
    @{Row=80; H=0.46885; I=0.01072},
    @{Row=81; H=0.00144; I=0.0033},
    @{Row=106; H=0.39635; I=0.02014},
    @{Row=107; H=0.01297; I=0.01053},
    @{Row=108; H=0.38018; I=0.0165},
    @{Row=109; H=0.00108; I=0.00293},
    @{Row=110; H=0.39815; I=0.02024},
    @{Row=111; H=0.01369; I=0.01054},
    @{Row=112; H=0.37939; I=0.01704},
    @{Row=113; H=0.0018; I=0.00441},
    @{Row=114; H=0.48161; I=0.01219},
    @{Row=115; H=0.00072; I=0.00244},
    @{Row=116; H=0.48166; I=0.01096},
    @{Row=117; H=0.00072; I=0.00244},
    @{Row=118; H=0.48178; I=0.01205},
    @{Row=119; H=0.00072; I=0.00244},
    @{Row=120; H=0.48178; I=0.01094},
    @{Row=121; H=0.00072; I=0.00244},
    @{Row=146; H=0.37023; I=0.02184},
    @{Row=147; H=0.01117; I=0.00817},
    @{Row=148; H=0.3619; I=0.02144},
    @{Row=149; H=0.00036; I=0.00177},
    @{Row=150; H=0.37331; I=0.02287},
    @{Row=151; H=0.01117; I=0.00817},
    @{Row=152; H=0.36389; I=0.02347},
    @{Row=153; H=0.00108; I=0.00293},
    @{Row=154; H=0.45944; I=0.01947},
    @{Row=155; H=0.00072; I=0.00244},
    @{Row=156; H=0.46455; I=0.01759},
    @{Row=157; H=0.00108; I=0.00293},
    @{Row=158; H=0.45963; I=0.01992},
    @{Row=159; H=0.00072; I=0.00244},
    @{Row=160; H=0.46369; I=0.01766},
    @{Row=161; H=0.00144; I=0.00293}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 8).Value = $u.H
    $ws.Cells.Item($u.Row, 9).Value = $u.I
}
